$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 165 (shifts existing rows 165-201 down to 166-202)
$ws.Rows.Item(165).Insert()

# Populate the newly inserted row 165 with the new record
$ws.Cells.Item(165, 1).Value = 11
$ws.Cells.Item(165, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(165, 3).Value = "Bíobío"
$ws.Cells.Item(165, 4).Value = 44995
$ws.Cells.Item(165, 5).Value = 8
$ws.Cells.Item(165, 6).Value = 100112032
$ws.Cells.Item(165, 7).Value = "Zapallo italiano"
$ws.Cells.Item(165, 8).Value = "Sin especificar"
$ws.Cells.Item(165, 9).Value = "Primera"
$ws.Cells.Item(165, 10).Value = 100
$ws.Cells.Item(165, 11).Value = 8000
$ws.Cells.Item(165, 12).Value = 8500
$ws.Cells.Item(165, 13).Value = 8250
$ws.Cells.Item(165, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(165, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(165, 16).Value = 165
$ws.Cells.Item(165, 17).Value = 50
$ws.Cells.Item(165, 18).Value = "Hortaliza"
